$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $value)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $ws.Range($cellRef).Value = "'" + $value
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

# Row 2
Set-CellText 'D2' '27.573.92'
Set-CellText 'E2' '  -1.72%  '

# Row 3
Set-CellText 'D3' '1.665.28'

# Row 4
Set-CellText 'E4' '  -0.01%  '

# Row 5
Set-CellText 'D5' '215.49'
Set-CellText 'E5' '  -1.53%  '

# Row 6
Set-CellText 'E6' '  -1.90%  '

# Row 7
Set-CellText 'E7' '  +0.00%  '

# Row 8
Set-CellText 'D8' '23.57'
Set-CellText 'E8' '  -2.33%  '

# Row 9
Set-CellText 'E9' '  -0.46%  '

# Row 10
Set-CellText 'D10' '0.0622'
Set-CellText 'E10' '  -1.74%  '

# Row 11
Set-CellText 'E11' '  -2.19%  '

# Row 12
Set-CellText 'D12' '1.900.60'
Set-CellText 'E12' '  -3.57%  '

# Row 13
Set-CellText 'D13' '1.682.78'
Set-CellText 'E13' '  -2.64%  '

# Row 14
Set-CellText 'E14' '  -2.74%  '

# Row 15
Set-CellText 'D15' '0.559'
Set-CellText 'E15' '  -0.81%  '

# Row 16
Set-CellText 'D16' '66.09'
Set-CellText 'E16' '  -2.28%  '

# Row 17
Set-CellText 'D17' '27.584.61'
Set-CellText 'E17' '  -1.56%  '

# Row 18
Set-CellText 'E18' '  -0.56%  '

# Row 19
Set-CellText 'E19' '  -3.43%  '

# Row 20
Set-CellText 'D20' '7.56'
Set-CellText 'E20' '  -3.99%  '

# Row 21
Set-CellText 'E21' '  +0.02%  '

# Row 22
Set-CellText 'E22' '  -3.26%  '

# Row 23
Set-CellText 'E23' '  -4.78%  '

# Row 24
Set-CellText 'E24' '  -3.94%  '

# Row 25
Set-CellText 'D25' '146.19'
Set-CellText 'E25' '  -1.93%  '

# Row 26
Set-CellText 'D26' '7.18'
Set-CellText 'E26' '  -4.55%  '

# Row 27
Set-CellText 'D27' '16.36'
Set-CellText 'E27' '  -2.31%  '

# Row 28
Set-CellText 'E28' '  +0.01%  '

# Row 29
Set-CellText 'E29' '  -2.43%  '

# Row 30
Set-CellText 'E30' '  +3.76%  '

# Row 31
Set-CellText 'D31' '0.0504'
Set-CellText 'E31' '  -1.22%  '

# Row 32
Set-CellText 'E32' '  -2.90%  '

# Row 33
Set-CellText 'D33' '1.478.93'
Set-CellText 'E33' '  -1.06%  '

# Row 34
Set-CellText 'D34' '3.10'
Set-CellText 'E34' '  -5.15%  '

# Row 35
Set-CellText 'E35' '  -5.58%  '

# Row 36
Set-CellText 'E36' '  -1.07%  '

# Row 37
Set-CellText 'D37' '0.932'
Set-CellText 'E37' '  -2.34%  '

# Row 38
Set-CellText 'B38' 'ImmutableX'
Set-CellText 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText 'D38' '0.573'
Set-CellText 'E38' '  -5.68%  '

# Row 39
Set-CellText 'B39' 'VeChain'
Set-CellText 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText 'D39' '0.0172'
Set-CellText 'E39' '  -1.86%  '

# Row 40
Set-CellText 'D40' '69.33'
Set-CellText 'E40' '  -2.15%  '

# Row 41
Set-CellText 'E41' '  -5.17%  '

# Row 42
Set-CellText 'E42' '  -0.05%  '

# Row 43
Set-CellText 'D43' '5.41'
Set-CellText 'E43' '  -7.33%  '

# Row 44
Set-CellText 'D44' '2.22'
Set-CellText 'E44' '  -3.70%  '

# Row 45
Set-CellText 'D45' '1.808.49'
Set-CellText 'E45' '  -3.50%  '

# Row 46
Set-CellText 'E46' '  -0.99%  '

# Row 47
Set-CellText 'E47' '  -3.05%  '

# Row 48
Set-CellText 'D48' '89.27'
Set-CellText 'E48' '  -2.10%  '

# Row 49
Set-CellText 'E49' '  -3.95%  '

# Row 50
Set-CellText 'E50' '  -2.60%  '

# Row 51
Set-CellText 'D51' '7.88'
Set-CellText 'E51' '  -3.43%  '
